# Update gh-pages to output generated at 456a3b4
#
# Adds the new "赣州·卡尼动漫展" (2024-09-15) event row to the "展览"
# (Exhibition) and "全部类型" (All types) sheets, and refreshes the
# "想去人数" (interested-count) figures that ticked up for several
# already-listed events across all affected sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $text) {
    # Force text storage so date-like strings ("2024-09-15") are not
    # auto-coerced into Excel date serials.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

function Insert-KanyiRow($ws, $rowNum, $idx) {
    # Insert a blank row at $rowNum (pushes everything below it down one),
    # copy column A's formatting down from the row above so the new index
    # cell keeps the same bold/bordered/centered style as its neighbours,
    # then populate the new row with the 赣州·卡尼动漫展 event data.
    $ws.Rows.Item($rowNum).Insert()

    $rowAbove = $rowNum - 1
    $addrAbove = "A" + $rowAbove
    $ws.Range($addrAbove).Copy()
    $addrA = "A" + $rowNum
    $ws.Range($addrA).PasteSpecial(-4122)
    $ws.Range($addrA).Value = $idx

    $addrB = "B" + $rowNum
    Set-TextCell $ws $addrB "2024-09-15"
    $addrC = "C" + $rowNum
    Set-TextCell $ws $addrC "赣州·卡尼动漫展"
    $addrD = "D" + $rowNum
    Set-TextCell $ws $addrD "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
    $addrE = "E" + $rowNum
    Set-TextCell $ws $addrE "2024.09.15 09:30-09.16 17:00"
    $addrF = "F" + $rowNum
    $ws.Range($addrF).Value = 23
    $addrG = "G" + $rowNum
    $ws.Range($addrG).Value = 45
    $addrH = "H" + $rowNum
    Set-TextCell $ws $addrH "https://show.bilibili.com/platform/detail.html?id=90642"
    $addrI = "I" + $rowNum
    Set-TextCell $ws $addrI "//i1.hdslb.com/bfs/openplatform/202408/VcJiaBPn1723530492504.jpeg"
}

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibition)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value = 5452
$ws1.Range("F8").Value = 898
$ws1.Range("F9").Value = 139
$ws1.Range("F10").Value = 2422

Insert-KanyiRow $ws1 12 11

# Rows 13 (old 12, unchanged) / 14 (old 13) / 15 (old 14) shifted down by
# the insert above; update the "想去人数" counts that changed.
$ws1.Range("F14").Value = 2271
$ws1.Range("F15").Value = 129

# ---------------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value = 96

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value = 5452
$ws4.Range("F6").Value = 96
$ws4.Range("F10").Value = 898
$ws4.Range("F11").Value = 139
$ws4.Range("F12").Value = 2422

Insert-KanyiRow $ws4 14 13

# Rows 15 (old 14) / 16 (old 15) unchanged counts; rows 17 (old 16) / 18
# (old 17) shifted down by the insert above and need updated counts.
$ws4.Range("F17").Value = 2271
$ws4.Range("F18").Value = 129
